$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10
$parts = @(
    "Please note: Every semester students fail to check first for the event ID = sl_bt_evt_system_external_signal_id and then check for the event value in evt->data.evt_system_external_signal.extsignals.",
    "",
    "Your state machine should call sl_bt_gatt_server_send_indication() if the appropriate conditions are met:",
    "Connection is open",
    "Client has enabled indications for the HTM indications",
    "There is no indication currently in-flight"
)
$text = [string]::Join($nl, $parts)
$rng = $ws.Range("G16")
$rng.Value = $text
$len = $text.Length
$c = $rng.Characters(1, $len)
$c.Font.Size = 11
$c.Font.Name = "Calibri"
$c.Font.ColorIndex = 8
Write-Output $rng.Value2
